$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("85:85").Insert()

$ws.Range("A85").Value = 4
$ws.Range("B85").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C85").Value = "Los Lagos"
$ws.Range("D85").Value = 44638
$ws.Range("E85").Value = 10
$ws.Range("F85").Value = "Fruta"
$ws.Range("G85").Value = 100102
$ws.Range("H85").Value = "Cítricos"
$ws.Range("I85").Value = 100102004
$ws.Range("J85").Value = "Mandarina"
$ws.Range("K85").Value = "Murcott"
$ws.Range("L85").Value = "Primera"
$ws.Range("M85").Value = 500
$ws.Range("N85").Value = 12500
$ws.Range("O85").Value = 13000
$ws.Range("P85").Value = 12750
$ws.Range("Q85").Value = "$/bandeja 10 kilos"
$ws.Range("R85").Value = "Región de O'Higgins"
$ws.Range("S85").Value = 1275
$ws.Range("T85").Value = 10

$ws.Range("D85").NumberFormat = $ws.Range("D86").NumberFormat
